$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = -10.97839999999999
$ws.Range("C10").Value = -12.2983
$ws.Range("C12").Value = -14.21920000000001
$ws.Range("D13").Value = -8.050600000000003
$ws.Range("C18").Value = -14.23720000000001
$ws.Range("C25").Value = -10.74709999999999
